$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Cells.Item(17, 8).Value = 812.9783
$ws.Cells.Item(17, 10).Value = 812.9783
$ws.Cells.Item(17, 12).Value = 2438.9349
$ws.Cells.Item(17, 14).Value = -2774.9349
# Row 18
$ws.Cells.Item(18, 8).Value = 1500
$ws.Cells.Item(18, 9).Value = 1500
$ws.Cells.Item(18, 11).Value = 1500
$ws.Cells.Item(18, 13).Value = -1216
# Row 32
$ws.Cells.Item(32, 8).Value = 4911.8945
$ws.Cells.Item(32, 10).Value = 4455.8
$ws.Cells.Item(32, 12).Value = 4455.8
$ws.Cells.Item(32, 14).Value = -5107.8
# Row 41
$ws.Cells.Item(41, 8).Value = 250.44444
$ws.Cells.Item(41, 9).Value = 244.28572
$ws.Cells.Item(41, 11).Value = 244.28572
$ws.Cells.Item(41, 13).Value = 195.71428
# Row 43
$ws.Cells.Item(43, 8).Value = 4973.125
$ws.Cells.Item(43, 10).Value = 5357
$ws.Cells.Item(43, 12).Value = 5357
$ws.Cells.Item(43, 14).Value = -5495
# Row 53
$ws.Cells.Item(53, 8).Value = 608.6667
$ws.Cells.Item(53, 9).Value = 620.8333
$ws.Cells.Item(53, 11).Value = 620.8333
$ws.Cells.Item(53, 13).Value = 16.16669999999999
# Row 74
$ws.Cells.Item(74, 8).Value = 10953.125
$ws.Cells.Item(74, 10).Value = 14600
$ws.Cells.Item(74, 12).Value = 14600
$ws.Cells.Item(74, 14).Value = -16472
# Row 77
$ws.Cells.Item(77, 8).Value = 10953.125
$ws.Cells.Item(77, 10).Value = 14600
$ws.Cells.Item(77, 12).Value = 73000
$ws.Cells.Item(77, 14).Value = -82360
# Row 92
$ws.Cells.Item(92, 8).Value = 205.54546
$ws.Cells.Item(92, 9).Value = 205.54546
$ws.Cells.Item(92, 11).Value = 205.54546
$ws.Cells.Item(92, 13).Value = 1042.45454
# Row 98
$ws.Cells.Item(98, 8).Value = 1912.2354
$ws.Cells.Item(98, 9).Value = 1977.9333
$ws.Cells.Item(98, 10).Value = 1419.5
$ws.Cells.Item(98, 11).Value = 1977.9333
$ws.Cells.Item(98, 12).Value = 1419.5
$ws.Cells.Item(98, 13).Value = -479.9332999999999
$ws.Cells.Item(98, 14).Value = -4415.5
# Row 122
$ws.Cells.Item(122, 8).Value = 1912.2354
$ws.Cells.Item(122, 9).Value = 1977.9333
$ws.Cells.Item(122, 10).Value = 1419.5
$ws.Cells.Item(122, 11).Value = 5933.7999
$ws.Cells.Item(122, 12).Value = 4258.5
$ws.Cells.Item(122, 13).Value = -3483.7999
$ws.Cells.Item(122, 14).Value = -9158.5
# Row 138
$ws.Cells.Item(138, 8).Value = 2710.25
$ws.Cells.Item(138, 10).Value = 0
$ws.Cells.Item(138, 12).Value = 0
$ws.Cells.Item(138, 14).ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 3554.5144
$ws.Cells.Item(32, 9).Value = 3457.8333
$ws.Cells.Item(32, 11).Value = 3457.8333
$ws.Cells.Item(32, 13).Value = -3170.8333
# Row 55
$ws.Cells.Item(55, 8).Value = 26687.8
$ws.Cells.Item(55, 10).Value = 34496.332
$ws.Cells.Item(55, 12).Value = 34496.332
$ws.Cells.Item(55, 14).Value = -35126.332
# Row 61
$ws.Cells.Item(61, 8).Value = 5771.5366
$ws.Cells.Item(61, 9).Value = 4700
$ws.Cells.Item(61, 11).Value = 4700
$ws.Cells.Item(61, 13).Value = -4488
# Row 81
$ws.Cells.Item(81, 8).Value = 60000
$ws.Cells.Item(81, 10).Value = 60000
$ws.Cells.Item(81, 12).Value = 60000
$ws.Cells.Item(81, 14).Value = -61996
# Row 84
$ws.Cells.Item(84, 8).Value = 60000
$ws.Cells.Item(84, 10).Value = 60000
$ws.Cells.Item(84, 12).Value = 180000
$ws.Cells.Item(84, 14).Value = -189984
# Row 136
$ws.Cells.Item(136, 8).Value = 5771.5366
$ws.Cells.Item(136, 9).Value = 4700
$ws.Cells.Item(136, 11).Value = 14100
$ws.Cells.Item(136, 13).Value = -11550

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Cells.Item(134, 8).Value = 3360.3157
$ws.Cells.Item(134, 9).Value = 3395.6667
$ws.Cells.Item(134, 11).Value = 10187.0001
$ws.Cells.Item(134, 13).Value = -7652.000100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 33
$ws.Cells.Item(33, 8).Value = 26000
$ws.Cells.Item(33, 10).Value = 26000
$ws.Cells.Item(33, 12).Value = 26000
$ws.Cells.Item(33, 14).Value = -26758
# Row 41
$ws.Cells.Item(41, 8).Value = 21958.143
$ws.Cells.Item(41, 10).Value = 21061.8
$ws.Cells.Item(41, 12).Value = 21061.8
$ws.Cells.Item(41, 14).Value = -21917.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Cells.Item(34, 8).Value = 298.33334
$ws.Cells.Item(34, 9).Value = 200
$ws.Cells.Item(34, 11).Value = 600
$ws.Cells.Item(34, 13).Value = -516
# Row 128
$ws.Cells.Item(128, 8).Value = 721831.5
$ws.Cells.Item(128, 9).Value = 721831.5
$ws.Cells.Item(128, 11).Value = 2165494.5
$ws.Cells.Item(128, 13).Value = -2160514.5
# Row 139
$ws.Cells.Item(139, 8).Value = 2640.05
$ws.Cells.Item(139, 9).Value = 2550.0557
$ws.Cells.Item(139, 11).Value = 7650.1671
$ws.Cells.Item(139, 13).Value = -2510.1671
# Row 140
$ws.Cells.Item(140, 8).Value = 1149.7667
$ws.Cells.Item(140, 9).Value = 696.2174
$ws.Cells.Item(140, 10).Value = 2640
$ws.Cells.Item(140, 11).Value = 2088.6522
$ws.Cells.Item(140, 12).Value = 7920
$ws.Cells.Item(140, 13).Value = 3091.3478
$ws.Cells.Item(140, 14).Value = -18280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3084.4
$ws.Cells.Item(80, 9).Value = 3160.4443
$ws.Cells.Item(80, 11).Value = 3160.4443
$ws.Cells.Item(80, 13).Value = -2162.4443
# Row 83
$ws.Cells.Item(83, 8).Value = 3084.4
$ws.Cells.Item(83, 9).Value = 3160.4443
$ws.Cells.Item(83, 11).Value = 15802.2215
$ws.Cells.Item(83, 13).Value = -10810.2215
# Row 132
$ws.Cells.Item(132, 8).Value = 3160.3667
$ws.Cells.Item(132, 9).Value = 2680.96
$ws.Cells.Item(132, 11).Value = 8042.88
$ws.Cells.Item(132, 13).Value = -5512.88

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 2125
$ws.Cells.Item(68, 9).Value = 2150
$ws.Cells.Item(68, 11).Value = 2150
$ws.Cells.Item(68, 13).Value = -1401
# Row 71
$ws.Cells.Item(71, 8).Value = 2125
$ws.Cells.Item(71, 9).Value = 2150
$ws.Cells.Item(71, 11).Value = 10750
$ws.Cells.Item(71, 13).Value = -7006
# Row 82
$ws.Cells.Item(82, 8).Value = 39999.5
$ws.Cells.Item(82, 9).Value = 39999.5
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 11).Value = 39999.5
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 13).Value = -39638.5
$ws.Cells.Item(82, 14).ClearContents()
# Row 85
$ws.Cells.Item(85, 8).Value = 39999.5
$ws.Cells.Item(85, 9).Value = 39999.5
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 39999.5
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 13).Value = -38751.5
$ws.Cells.Item(85, 14).ClearContents()
# Row 122
$ws.Cells.Item(122, 8).Value = 3162.35
$ws.Cells.Item(122, 9).Value = 3091.389
$ws.Cells.Item(122, 11).Value = 9274.167000000001
$ws.Cells.Item(122, 13).Value = -6824.167000000001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Cells.Item(41, 8).Value = 10956.692
$ws.Cells.Item(41, 9).Value = 9000
$ws.Cells.Item(41, 10).Value = 11119.75
$ws.Cells.Item(41, 11).Value = 9000
$ws.Cells.Item(41, 12).Value = 11119.75
$ws.Cells.Item(41, 13).Value = -8610
$ws.Cells.Item(41, 14).Value = -11899.75
# Row 62
$ws.Cells.Item(62, 8).Value = 4625.778
$ws.Cells.Item(62, 9).Value = 3954
$ws.Cells.Item(62, 10).Value = 10000
$ws.Cells.Item(62, 11).Value = 3954
$ws.Cells.Item(62, 12).Value = 10000
$ws.Cells.Item(62, 13).Value = -3330
$ws.Cells.Item(62, 14).Value = -11248
# Row 65
$ws.Cells.Item(65, 8).Value = 4625.778
$ws.Cells.Item(65, 9).Value = 3954
$ws.Cells.Item(65, 10).Value = 10000
$ws.Cells.Item(65, 11).Value = 19770
$ws.Cells.Item(65, 12).Value = 50000
$ws.Cells.Item(65, 13).Value = -16650
$ws.Cells.Item(65, 14).Value = -56240
# Row 81
$ws.Cells.Item(81, 8).Value = 749.6667
$ws.Cells.Item(81, 9).Value = 749.6667
$ws.Cells.Item(81, 11).Value = 1499.3334
$ws.Cells.Item(81, 13).Value = -438.3334
# Row 84
$ws.Cells.Item(84, 8).Value = 749.6667
$ws.Cells.Item(84, 9).Value = 749.6667
$ws.Cells.Item(84, 11).Value = 7496.666999999999
$ws.Cells.Item(84, 13).Value = -2192.666999999999
# Row 100
$ws.Cells.Item(100, 8).Value = 1257.8462
$ws.Cells.Item(100, 9).Value = 1255.8
$ws.Cells.Item(100, 11).Value = 2511.6
$ws.Cells.Item(100, 13).Value = -1970.6
# Row 107
$ws.Cells.Item(107, 8).Value = 555.06665
$ws.Cells.Item(107, 9).Value = 511
$ws.Cells.Item(107, 10).Value = 841.5
$ws.Cells.Item(107, 11).Value = 1533
$ws.Cells.Item(107, 12).Value = 2524.5
$ws.Cells.Item(107, 13).Value = 387
$ws.Cells.Item(107, 14).Value = -6364.5
# Row 126
$ws.Cells.Item(126, 8).Value = 3667.6667
$ws.Cells.Item(126, 9).Value = 3667.6667
$ws.Cells.Item(126, 11).Value = 11003.0001
$ws.Cells.Item(126, 13).Value = -8533.000100000001
# Row 132
$ws.Cells.Item(132, 8).Value = 2392.5925
$ws.Cells.Item(132, 9).Value = 2476.077
$ws.Cells.Item(132, 10).Value = 222
$ws.Cells.Item(132, 11).Value = 7428.231000000001
$ws.Cells.Item(132, 12).Value = 666
$ws.Cells.Item(132, 13).Value = -4898.231000000001
$ws.Cells.Item(132, 14).Value = -5726
Write-Host "All updates applied"
